$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 64
$ws.Range("H64").Value = 4871.2856
$ws.Range("I64").Value = 4766.5
$ws.Range("J64").Value = 5500
$ws.Range("K64").Value = 4766.5
$ws.Range("L64").Value = 5500
$ws.Range("M64").Value = -4518.5
$ws.Range("N64").Value = -5996
# Row 67
$ws.Range("H67").Value = 4871.2856
$ws.Range("I67").Value = 4766.5
$ws.Range("J67").Value = 5500
$ws.Range("K67").Value = 4766.5
$ws.Range("L67").Value = 5500
$ws.Range("M67").Value = -3908.5
$ws.Range("N67").Value = -7216
# Row 69
$ws.Range("H69").Value = 6593.3335
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 6593.3335
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 19780.0005
$ws.Range("N69").Value = -21528.0005
$ws.Range("M69").ClearContents()
# Row 72
$ws.Range("H72").Value = 6593.3335
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 6593.3335
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 59340.0015
$ws.Range("N72").Value = -68076.0015
$ws.Range("M72").ClearContents()
# Row 76
$ws.Range("H76").Value = 68184560
$ws.Range("I76").Value = 75002696
$ws.Range("J76").Value = 3250
$ws.Range("K76").Value = 75002696
$ws.Range("L76").Value = 3250
$ws.Range("M76").Value = -75002381
$ws.Range("N76").Value = -3880
# Row 79
$ws.Range("H79").Value = 68184560
$ws.Range("I79").Value = 75002696
$ws.Range("J79").Value = 3250
$ws.Range("K79").Value = 75002696
$ws.Range("L79").Value = 3250
$ws.Range("M79").Value = -75001604
$ws.Range("N79").Value = -5434
# Row 137
$ws.Range("H137").Value = 846.0732
$ws.Range("I137").Value = 790
$ws.Range("J137").Value = 1019.9
$ws.Range("K137").Value = 2370
$ws.Range("L137").Value = 3059.7
$ws.Range("M137").Value = 180
$ws.Range("N137").Value = -8159.7

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 2857
$ws.Range("I61").Value = 3230.2856
$ws.Range("J61").Value = 2421.5
$ws.Range("K61").Value = 3230.2856
$ws.Range("L61").Value = 2421.5
$ws.Range("M61").Value = -3018.2856
$ws.Range("N61").Value = -2845.5
# Row 63
$ws.Range("H63").Value = 2053.4707
$ws.Range("I63").Value = 2053.4707
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 2053.4707
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -1367.4707
$ws.Range("N63").ClearContents()
# Row 66
$ws.Range("H66").Value = 2053.4707
$ws.Range("I66").Value = 2053.4707
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 10267.3535
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -6835.353499999999
$ws.Range("N66").ClearContents()
# Row 88
$ws.Range("H88").Value = 2485.8125
$ws.Range("I88").Value = 1982.375
$ws.Range("J88").Value = 2989.25
$ws.Range("K88").Value = 1982.375
$ws.Range("L88").Value = 2989.25
$ws.Range("M88").Value = -1576.375
$ws.Range("N88").Value = -3801.25
# Row 91
$ws.Range("H91").Value = 2485.8125
$ws.Range("I91").Value = 1982.375
$ws.Range("J91").Value = 2989.25
$ws.Range("K91").Value = 1982.375
$ws.Range("L91").Value = 2989.25
$ws.Range("M91").Value = -578.375
$ws.Range("N91").Value = -5797.25
# Row 132
$ws.Range("H132").Value = 1471.1968
$ws.Range("I132").Value = 1202.561
$ws.Range("J132").Value = 2021.9
$ws.Range("K132").Value = 3607.683
$ws.Range("L132").Value = 6065.700000000001
$ws.Range("M132").Value = -1077.683
$ws.Range("N132").Value = -11125.7
# Row 136
$ws.Range("H136").Value = 2857
$ws.Range("I136").Value = 3230.2856
$ws.Range("J136").Value = 2421.5
$ws.Range("K136").Value = 9690.856800000001
$ws.Range("L136").Value = 7264.5
$ws.Range("M136").Value = -7140.856800000001
$ws.Range("N136").Value = -12364.5

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 6011.1
$ws.Range("I105").Value = 6467.8887
$ws.Range("J105").Value = 1900
$ws.Range("K105").Value = 6467.8887
$ws.Range("L105").Value = 1900
$ws.Range("M105").Value = -4720.8887
$ws.Range("N105").Value = -5394

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 31486
$ws.Range("I31").Value = 3154.16
$ws.Range("J31").Value = 102315.6
$ws.Range("K31").Value = 3154.16
$ws.Range("L31").Value = 102315.6
$ws.Range("M31").Value = -2859.16
$ws.Range("N31").Value = -102905.6
# Row 34
$ws.Range("H34").Value = 31486
$ws.Range("I34").Value = 3154.16
$ws.Range("J34").Value = 102315.6
$ws.Range("K34").Value = 3154.16
$ws.Range("L34").Value = 102315.6
$ws.Range("M34").Value = -2952.16
$ws.Range("N34").Value = -102719.6
# Row 62
$ws.Range("H62").Value = 4364.077
$ws.Range("I62").Value = 3519.125
$ws.Range("J62").Value = 5716
$ws.Range("K62").Value = 3519.125
$ws.Range("L62").Value = 5716
$ws.Range("M62").Value = -2895.125
$ws.Range("N62").Value = -6964
# Row 65
$ws.Range("H65").Value = 4364.077
$ws.Range("I65").Value = 3519.125
$ws.Range("J65").Value = 5716
$ws.Range("K65").Value = 17595.625
$ws.Range("L65").Value = 28580
$ws.Range("M65").Value = -14475.625
$ws.Range("N65").Value = -34820
# Row 134
$ws.Range("H134").Value = 41668036
$ws.Range("I134").Value = 1495.1818
$ws.Range("J134").Value = 500000000
$ws.Range("K134").Value = 4485.5454
$ws.Range("L134").Value = 1500000000
$ws.Range("M134").Value = -1950.5454
$ws.Range("N134").Value = -1500005070

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 30898.4
$ws.Range("I70").Value = 53571
$ws.Range("J70").Value = 4986.857
$ws.Range("K70").Value = 53571
$ws.Range("L70").Value = 4986.857
$ws.Range("M70").Value = -53301
$ws.Range("N70").Value = -5526.857
# Row 73
$ws.Range("H73").Value = 30898.4
$ws.Range("I73").Value = 53571
$ws.Range("J73").Value = 4986.857
$ws.Range("K73").Value = 53571
$ws.Range("L73").Value = 4986.857
$ws.Range("M73").Value = -52635
$ws.Range("N73").Value = -6858.857
# Row 80
$ws.Range("H80").Value = 3524.318
$ws.Range("I80").Value = 3719.1177
$ws.Range("J80").Value = 2862
$ws.Range("K80").Value = 3719.1177
$ws.Range("L80").Value = 2862
$ws.Range("M80").Value = -2721.1177
$ws.Range("N80").Value = -4858
# Row 83
$ws.Range("H83").Value = 3524.318
$ws.Range("I83").Value = 3719.1177
$ws.Range("J83").Value = 2862
$ws.Range("K83").Value = 18595.5885
$ws.Range("L83").Value = 14310
$ws.Range("M83").Value = -13603.5885
$ws.Range("N83").Value = -24294

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 132
$ws.Range("H132").Value = 3074.3777
$ws.Range("I132").Value = 2056.7727
$ws.Range("J132").Value = 4047.739
$ws.Range("K132").Value = 6170.3181
$ws.Range("L132").Value = 12143.217
$ws.Range("M132").Value = -3640.3181
$ws.Range("N132").Value = -17203.217

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 41
$ws.Range("H41").Value = 13164.857
$ws.Range("I41").Value = 9000
$ws.Range("J41").Value = 13859
$ws.Range("K41").Value = 9000
$ws.Range("L41").Value = 13859
$ws.Range("M41").Value = -8610
$ws.Range("N41").Value = -14639
# Row 46
$ws.Range("H46").Value = 45799.855
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 45799.855
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 45799.855
$ws.Range("N46").Value = -46261.855
# Row 134
$ws.Range("H134").Value = 45799.855
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 45799.855
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 137399.565
$ws.Range("N134").Value = -142469.565
# Row 135
$ws.Range("H135").Value = 75212.78
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 75212.78
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 75212.78
$ws.Range("N135").Value = -85352.78
